$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Title Page")

# Correct the wrong document name on the Title Page: the review sheet
# referenced "CommLib API Android" instead of "Technical Design CommLib
# Android" in the document/product header cell (A6, merged A6:H6).
$ws1.Range("A6").Value = "Document: CML000009 Technical Design CommLib Android`nProduct/Platform: CommLib Android/CDP2 Platform"

# Move the on-screen selection to the corrected cell/row.
$ws1.Activate()
$ws1.Range("A6:H6").Select()
